# Update Leve profit-calculation market price snapshots across all crafting-class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the latest market data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4307.4473
$ws.Range("J32").Value = 7490.857
$ws.Range("L32").Value = 7490.857
$ws.Range("N32").Value = -8064.857
$ws.Range("H61").Value = 1195.0834
$ws.Range("I61").Value = 1144.6818
$ws.Range("K61").Value = 1144.6818
$ws.Range("M61").Value = -932.6818000000001
$ws.Range("H109").Value = 39750
$ws.Range("J109").Value = 39750
$ws.Range("L109").Value = 39750
$ws.Range("N109").Value = -42524
$ws.Range("H113").Value = 30358.2
$ws.Range("J113").Value = 30358.2
$ws.Range("L113").Value = 30358.2
$ws.Range("N113").Value = -39036.2
$ws.Range("H136").Value = 1195.0834
$ws.Range("I136").Value = 1144.6818
$ws.Range("K136").Value = 3434.0454
$ws.Range("M136").Value = -884.0454

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17217.182
$ws.Range("I82").Value = 13628.5
$ws.Range("J82").Value = 18014.666
$ws.Range("K82").Value = 13628.5
$ws.Range("L82").Value = 18014.666
$ws.Range("M82").Value = -13245.5
$ws.Range("N82").Value = -18780.666
$ws.Range("H85").Value = 17217.182
$ws.Range("I85").Value = 13628.5
$ws.Range("J85").Value = 18014.666
$ws.Range("K85").Value = 13628.5
$ws.Range("L85").Value = 18014.666
$ws.Range("M85").Value = -12302.5
$ws.Range("N85").Value = -20666.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1322.6666
$ws.Range("I16").Value = 1292.3889
$ws.Range("K16").Value = 1292.3889
$ws.Range("M16").Value = -1005.3889
$ws.Range("H22").Value = 284
$ws.Range("I22").Value = 271.45456
$ws.Range("K22").Value = 271.45456
$ws.Range("M22").Value = 78.54543999999999
$ws.Range("H31").Value = 2142.3667
$ws.Range("I31").Value = 2130.0344
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 2130.0344
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1835.0344
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 2142.3667
$ws.Range("I34").Value = 2130.0344
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 2130.0344
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1928.0344
$ws.Range("N34").Value = -2904
$ws.Range("H58").Value = 4069.2122
$ws.Range("I58").Value = 673.3333
$ws.Range("J58").Value = 13124.889
$ws.Range("K58").Value = 673.3333
$ws.Range("L58").Value = 13124.889
$ws.Range("M58").Value = -470.3333
$ws.Range("N58").Value = -13530.889
$ws.Range("H113").Value = 1322.6666
$ws.Range("I113").Value = 1292.3889
$ws.Range("K113").Value = 1292.3889
$ws.Range("M113").Value = 877.6111000000001
$ws.Range("H134").Value = 1799.1364
$ws.Range("I134").Value = 1284.05
$ws.Range("K134").Value = 3852.15
$ws.Range("M134").Value = -1317.15
$ws.Range("H136").Value = 4069.2122
$ws.Range("I136").Value = 673.3333
$ws.Range("J136").Value = 13124.889
$ws.Range("K136").Value = 2019.9999
$ws.Range("L136").Value = 39374.667
$ws.Range("M136").Value = 530.0001
$ws.Range("N136").Value = -44474.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 138.4
$ws.Range("I23").Value = 97.5
$ws.Range("J23").Value = 302
$ws.Range("K23").Value = 292.5
$ws.Range("L23").Value = 906
$ws.Range("M23").Value = -57.5
$ws.Range("N23").Value = -1376
$ws.Range("H109").Value = 2544.3333
$ws.Range("I109").Value = 699.5
$ws.Range("J109").Value = 3071.4285
$ws.Range("K109").Value = 2098.5
$ws.Range("L109").Value = 9214.2855
$ws.Range("M109").Value = -1058.5
$ws.Range("N109").Value = -11294.2855
$ws.Range("H117").Value = 2224.05
$ws.Range("J117").Value = 2317.9443
$ws.Range("L117").Value = 6953.8329
$ws.Range("N117").Value = -13837.8329
$ws.Range("H129").Value = 1819.3334
$ws.Range("I129").Value = 1047.9
$ws.Range("J129").Value = 2154.739
$ws.Range("K129").Value = 3143.7
$ws.Range("L129").Value = 6464.217000000001
$ws.Range("M129").Value = 1856.3
$ws.Range("N129").Value = -16464.217
$ws.Range("H131").Value = 21446.156
$ws.Range("I131").Value = 92179.91
$ws.Range("J131").Value = 1994.375
$ws.Range("K131").Value = 276539.73
$ws.Range("L131").Value = 5983.125
$ws.Range("M131").Value = -271499.73
$ws.Range("N131").Value = -16063.125
$ws.Range("H132").Value = 3333901.8
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 5000502.5
$ws.Range("K132").Value = 6300
$ws.Range("L132").Value = 45004522.5
$ws.Range("M132").Value = -3770
$ws.Range("N132").Value = -45009582.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4377.6665
$ws.Range("I70").Value = 4087.88
$ws.Range("K70").Value = 4087.88
$ws.Range("M70").Value = -3817.88
$ws.Range("H73").Value = 4377.6665
$ws.Range("I73").Value = 4087.88
$ws.Range("K73").Value = 4087.88
$ws.Range("M73").Value = -3151.88
$ws.Range("H123").Value = 25986
$ws.Range("J123").Value = 25986
$ws.Range("L123").Value = 25986
$ws.Range("N123").Value = -30886

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 749.8
$ws.Range("I22").Value = 685.4286
$ws.Range("K22").Value = 685.4286
$ws.Range("M22").Value = -390.4286
$ws.Range("H27").Value = 749.8
$ws.Range("I27").Value = 685.4286
$ws.Range("K27").Value = 685.4286
$ws.Range("M27").Value = -578.4286
$ws.Range("H68").Value = 2873
$ws.Range("I68").Value = 3534
$ws.Range("J68").Value = 2625.125
$ws.Range("K68").Value = 3534
$ws.Range("L68").Value = 2625.125
$ws.Range("M68").Value = -2785
$ws.Range("N68").Value = -4123.125
$ws.Range("H71").Value = 2873
$ws.Range("I71").Value = 3534
$ws.Range("J71").Value = 2625.125
$ws.Range("K71").Value = 17670
$ws.Range("L71").Value = 13125.625
$ws.Range("M71").Value = -13926
$ws.Range("N71").Value = -20613.625
$ws.Range("H101").Value = 17120.666
$ws.Range("J101").Value = 17120.666
$ws.Range("L101").Value = 17120.666
$ws.Range("N101").Value = -23610.666
$ws.Range("H136").Value = 2285.8975
$ws.Range("J136").Value = 6786.4287
$ws.Range("L136").Value = 20359.2861
$ws.Range("N136").Value = -25459.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5222.222
$ws.Range("I62").Value = 2800
$ws.Range("J62").Value = 7160
$ws.Range("K62").Value = 2800
$ws.Range("L62").Value = 7160
$ws.Range("M62").Value = -2176
$ws.Range("N62").Value = -8408
$ws.Range("H65").Value = 5222.222
$ws.Range("I65").Value = 2800
$ws.Range("J65").Value = 7160
$ws.Range("K65").Value = 14000
$ws.Range("L65").Value = 35800
$ws.Range("M65").Value = -10880
$ws.Range("N65").Value = -42040
$ws.Range("H103").Value = 38200.5
$ws.Range("J103").Value = 38200.5
$ws.Range("L103").Value = 38200.5
$ws.Range("N103").Value = -40544.5
